$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Relevant Parts")

# Add new shared string / cell value: digikey link for the LED part (row 10, PWR LED)
$ws.Range("B10").Value = "https://www.digikey.com/en/products/detail/american-opto-plus-led/L171L-GC/12325425"

# Update selection to match the saved cursor position
$ws.Range("F17").Select()
